$wb = $excel.ActiveWorkbook

# --- GroupTable: drop quotes from group_name values, widen column B ---
$ws2 = $wb.Worksheets.Item("GroupTable")
$ws2.Range("B2").Value = "3TeamExtreme"
$ws2.Range("B3").Value = "Lab42"
$ws2.Range("B4").Value = "Misfits"
$ws2.Columns.Item(2).ColumnWidth = 14.166666666666666

# --- UserTable: drop quotes from user_name / site values ---
# A2/A3 keep their pre-existing "quote prefix" text format, so re-enter them
# with a leading apostrophe (a pure formatting marker Excel strips from the
# stored text) to preserve that cell style.
$ws1 = $wb.Worksheets.Item("UserTable")
$ws1.Range("A2").Value = "'Bob"
$ws1.Range("A3").Value = "'Mary"
$ws1.Range("A4").Value = "John"
$ws1.Range("A5").Value = "Frank"
$ws1.Range("A6").Value = "Sally"
$ws1.Range("B2").Value = "Lake Mary"
$ws1.Range("B4").Value = "Lake Mary"
$ws1.Range("B5").Value = "Sanford"

# --- RoleTable: drop quotes and capitalize role_name values ---
$ws4 = $wb.Worksheets.Item("RoleTable")
$ws4.Range("B2").Value = "Chemist"
$ws4.Range("B3").Value = "Biologist"
$ws4.Range("B4").Value = "Analyst"

# site values containing "Winter Park" are re-typed last (matches authoring order)
$ws1.Range("B3").Value = "Winter Park"
$ws1.Range("B6").Value = "Winter Park"

# --- GroupJunctionTable: drop quotes from username values ---
$ws3 = $wb.Worksheets.Item("GroupJunctionTable")
$ws3.Range("A2").Value = "Bob"
$ws3.Range("A3").Value = "Mary"
$ws3.Range("A4").Value = "Mary"
$ws3.Range("A5").Value = "Mary"
$ws3.Range("A6").Value = "John"
$ws3.Range("A7").Value = "Frank"
$ws3.Range("A8").Value = "Frank"
$ws3.Range("A9").Value = "Sally"

# --- RoleJunctionTable: drop quotes from username values ---
$ws5 = $wb.Worksheets.Item("RoleJunctionTable")
$ws5.Range("A2").Value = "Bob"
$ws5.Range("A3").Value = "Mary"
$ws5.Range("A4").Value = "Mary"
$ws5.Range("A5").Value = "John"
$ws5.Range("A6").Value = "John"
$ws5.Range("A7").Value = "Frank"
$ws5.Range("A8").Value = "Sally"
$ws5.Range("A9").Value = "Sally"

# --- Restore per-sheet selections (view-only changes from the diff) ---
$ws2.Activate()
$ws2.Range("C7").Select()

$ws3.Activate()
$ws3.Range("A10").Select()

$ws4.Activate()
$ws4.Range("B5").Select()

$ws5.Activate()
$ws5.Range("A9").Select()

$ws1.Activate()
$ws1.Range("C3").Select()
